$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1317.5625
$ws.Range("I107").Value = 831.36365
$ws.Range("K107").Value = 831.36365
$ws.Range("M107").Value = 1088.63635
$ws.Range("H137").Value = 1759.8462
$ws.Range("I137").Value = 1510.6842
$ws.Range("J137").Value = 2436.1428
$ws.Range("K137").Value = 4532.0526
$ws.Range("L137").Value = 7308.428400000001
$ws.Range("M137").Value = -1982.0526
$ws.Range("N137").Value = -12408.4284
$ws.Range("H138").Value = 1969.804
$ws.Range("I138").Value = 1464.762
$ws.Range("J138").Value = 2323.3333
$ws.Range("K138").Value = 4394.286
$ws.Range("L138").Value = 6969.999899999999
$ws.Range("M138").Value = 745.7139999999999
$ws.Range("N138").Value = -17249.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5742.5
$ws.Range("I32").Value = 3522.9788
$ws.Range("K32").Value = 3522.9788
$ws.Range("M32").Value = -3235.9788
$ws.Range("H74").Value = 1239.1428
$ws.Range("I74").Value = 565.1739
$ws.Range("K74").Value = 565.1739
$ws.Range("M74").Value = 308.8261
$ws.Range("H77").Value = 1239.1428
$ws.Range("I77").Value = 565.1739
$ws.Range("K77").Value = 2825.8695
$ws.Range("M77").Value = 1542.1305
$ws.Range("H110").Value = 325
$ws.Range("I110").Value = 325
$ws.Range("K110").Value = 325
$ws.Range("M110").Value = 1720
$ws.Range("H122").Value = 998.2692
$ws.Range("I122").Value = 807.2
$ws.Range("K122").Value = 2421.6
$ws.Range("M122").Value = 28.39999999999964
$ws.Range("H132").Value = 1789.5555
$ws.Range("I132").Value = 1371.762
$ws.Range("K132").Value = 4115.286
$ws.Range("M132").Value = -1585.286

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2428.5454
$ws.Range("I105").Value = 1857.1111
$ws.Range("K105").Value = 1857.1111
$ws.Range("M105").Value = -110.1111000000001
$ws.Range("H107").Value = 856.2353
$ws.Range("I107").Value = 597.6667
$ws.Range("J107").Value = 1147.125
$ws.Range("K107").Value = 597.6667
$ws.Range("L107").Value = 1147.125
$ws.Range("M107").Value = 1322.3333
$ws.Range("N107").Value = -4987.125
$ws.Range("H134").Value = 6604.8184
$ws.Range("I134").Value = 8643.214
$ws.Range("J134").Value = 3037.625
$ws.Range("K134").Value = 25929.642
$ws.Range("L134").Value = 9112.875
$ws.Range("M134").Value = -23394.642
$ws.Range("N134").Value = -14182.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 879.9
$ws.Range("I105").Value = 904
$ws.Range("J105").Value = 823.6667
$ws.Range("K105").Value = 904
$ws.Range("L105").Value = 823.6667
$ws.Range("M105").Value = 843
$ws.Range("N105").Value = -4317.6667
$ws.Range("H107").Value = 1012.5455
$ws.Range("I107").Value = 727.8
$ws.Range("J107").Value = 1249.8334
$ws.Range("K107").Value = 727.8
$ws.Range("L107").Value = 1249.8334
$ws.Range("M107").Value = 1192.2
$ws.Range("N107").Value = -5089.8334
$ws.Range("H134").Value = 3000
$ws.Range("J134").Value = 3000
$ws.Range("L134").Value = 9000
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1212.875
$ws.Range("I97").Value = 1279
$ws.Range("K97").Value = 1279
$ws.Range("M97").Value = -783
$ws.Range("H122").Value = 1461.35
$ws.Range("I122").Value = 1538.5834
$ws.Range("J122").Value = 1345.5
$ws.Range("K122").Value = 4615.7502
$ws.Range("L122").Value = 4036.5
$ws.Range("M122").Value = -2165.7502
$ws.Range("N122").Value = -8936.5
$ws.Range("H132").Value = 3849423
$ws.Range("J132").Value = 4449.3335
$ws.Range("L132").Value = 13348.0005
$ws.Range("N132").Value = -18408.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 12998.308
$ws.Range("I40").Value = 12280.363
$ws.Range("K40").Value = 12280.363
$ws.Range("M40").Value = -12144.363
$ws.Range("H61").Value = 2325.7058
$ws.Range("I61").Value = 2103.7
$ws.Range("J61").Value = 2642.8572
$ws.Range("K61").Value = 2103.7
$ws.Range("L61").Value = 2642.8572
$ws.Range("M61").Value = -1901.7
$ws.Range("N61").Value = -3046.8572
$ws.Range("H98").Value = 30000
$ws.Range("J98").Value = 30000
$ws.Range("L98").Value = 30000
$ws.Range("N98").Value = -35990
$ws.Range("H113").Value = 2325.7058
$ws.Range("I113").Value = 2103.7
$ws.Range("J113").Value = 2642.8572
$ws.Range("K113").Value = 2103.7
$ws.Range("L113").Value = 2642.8572
$ws.Range("M113").Value = 66.30000000000018
$ws.Range("N113").Value = -6982.8572
$ws.Range("H122").Value = 5021.24
$ws.Range("I122").Value = 4418.1113
$ws.Range("K122").Value = 13254.3339
$ws.Range("M122").Value = -10804.3339

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 50005000
$ws.Range("I62").Value = 100000000
$ws.Range("J62").Value = 10000
$ws.Range("K62").Value = 100000000
$ws.Range("L62").Value = 10000
$ws.Range("M62").Value = -99999376
$ws.Range("N62").Value = -11248
$ws.Range("H65").Value = 50005000
$ws.Range("I65").Value = 100000000
$ws.Range("J65").Value = 10000
$ws.Range("K65").Value = 500000000
$ws.Range("L65").Value = 50000
$ws.Range("M65").Value = -499996880
$ws.Range("N65").Value = -56240
$ws.Range("H68").Value = 85000
$ws.Range("J68").Value = 85000
$ws.Range("L68").Value = 85000
$ws.Range("N68").Value = -86622
$ws.Range("H71").Value = 85000
$ws.Range("J71").Value = 85000
$ws.Range("L71").Value = 255000
$ws.Range("N71").Value = -263112
$ws.Range("H107").Value = 940.2308
$ws.Range("I107").Value = 469.44446
$ws.Range("J107").Value = 1999.5
$ws.Range("K107").Value = 1408.33338
$ws.Range("L107").Value = 5998.5
$ws.Range("M107").Value = 511.66662
$ws.Range("N107").Value = -9838.5
$ws.Range("H113").Value = 441.7931
$ws.Range("I113").Value = 312.11765
$ws.Range("J113").Value = 625.5
$ws.Range("K113").Value = 936.3529500000001
$ws.Range("L113").Value = 1876.5
$ws.Range("M113").Value = 1233.64705
$ws.Range("N113").Value = -6216.5
$ws.Range("H122").Value = 30964
$ws.Range("I122").Value = 37899.145
$ws.Range("J122").Value = 1836.4
$ws.Range("K122").Value = 113697.435
$ws.Range("L122").Value = 5509.200000000001
$ws.Range("M122").Value = -111247.435
$ws.Range("N122").Value = -10409.2
$ws.Range("H136").Value = 17922432
$ws.Range("I136").Value = 25253442
$ws.Range("J136").Value = 2187.6667
$ws.Range("K136").Value = 75760326
$ws.Range("L136").Value = 6563.000100000001
$ws.Range("M136").Value = -75757776
$ws.Range("N136").Value = -11663.0001
